$wb = $excel.ActiveWorkbook

$sheetNames = @("pre", "standard")

# New data (by sheet name) taken from the target OOXML.
$data = @{
    "pre" = @{
        Header = @("China", "EU", "India")
        Rows = @(
            @{ Label = "Y";        Values = @(-9.932487581726193, 24.98504778026656, -15.05256019854037) },
            @{ Label = "Y_Import"; Values = @(6.97521704121815, -14.14170446634398, 7.166487425125831) },
            @{ Label = "r_c";      Values = @(-5.032224880158286, 10.84554644021556, -5.813321560057277) },
            @{ Label = "r_p";      Values = @(-3.436828665498592, 8.7831677131497, -5.346339047651106) }
        )
    }
    "standard" = @{
        Header = @("China", "EU", "India")
        Rows = @(
            @{ Label = "Y";        Values = @(-1.14860886739979, 1.238235081695834, -1.302697326032186) },
            @{ Label = "Y_Import"; Values = @(1.5959637597411, -1.551424851074756, 1.509890094087287) },
            @{ Label = "r_c";      Values = @(-0.3531651820238457, 0.2301166437257283, -0.13315276200683) },
            @{ Label = "r_p";      Values = @(-0.09418971031746368, 0.08307312565319393, -0.07404000604827089) }
        )
    }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cfg = $data[$sheetName]

    # Drop column E entirely (Y_Import/r_p's old "r_p" header column no longer exists)
    # and drop A1 ("Region" label), since the new layout has no row/column header corner cell.
    $ws.Columns.Item(5).Clear()
    $ws.Range("A1").Clear()

    # Row 1: reuse existing styled cells B1:D1 (keep their bold/border/center style),
    # just replace the text with the new column headers.
    $ws.Range("B1").Value = $cfg.Header[0]
    $ws.Range("C1").Value = $cfg.Header[1]
    $ws.Range("D1").Value = $cfg.Header[2]

    # Column A labels: reuse the three existing styled label cells (A2:A4) in place,
    # then create a new styled A5 cell (copying the format from A4) for the 4th label.
    $ws.Range("A2").Value = $cfg.Rows[0].Label
    $ws.Range("A3").Value = $cfg.Rows[1].Label
    $ws.Range("A4").Value = $cfg.Rows[2].Label

    $ws.Range("A4").Copy()
    $ws.Range("A5").PasteSpecial(-4122)
    $ws.Range("A5").Value = $cfg.Rows[3].Label

    # Numeric data, rows 2-5, columns B:D
    for ($i = 0; $i -lt 4; $i++) {
        $r = 2 + $i
        $vals = $cfg.Rows[$i].Values
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
    }
}
